$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '62.871.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.01%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.012.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -2.41%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.01%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '557.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +0.40%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '153.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -3.58%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.03%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -2.11%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '3.014.14'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -2.34%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -0.61%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -4.06%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -1.55%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '3.534.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -2.69%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D15').Value = "'" + '62.963.92'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +0.05%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '23.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -1.19%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '3.014.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -2.14%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '0.0000149'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -0.45%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '395.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +1.10%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -0.07%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '11.89'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -2.52%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '6.63'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -3.77%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -0.22%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '65.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -2.31%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -0.25%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -3.91%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '0.0₃0971'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -1.63%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '8.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +2.03%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -0.63%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +0.06%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '1.74'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -0.68%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -0.38%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '159.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +5.88%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '4.69'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -0.46%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -1.29%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +1.58%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '1.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +0.44%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '2.507.87'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -5.87%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -2.05%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '37.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -1.27%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '22.51'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -1.59%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '3.91'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -2.00%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '0.665'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -3.06%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -0.16%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '0.998'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -0.16%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'" + 'VeChain'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'" + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'" + '0.0246'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -1.75%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'" + 'RenderToken'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'" + '5.03'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -6.16%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '19.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -2.62%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'" + 'WhiteBITCoin'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'" + 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'" + '10.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +0.47%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'" + 'Stellar'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'" + '0.0946'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -2.03%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '263.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -5.51%  '
$ws.Range('E51').Style = 'Normal'
